$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Res_CC"
$ws.Range("E1").Value = "Res_II_1"
$ws.Range("G1").Value = "Res_II_2"

$ws.Range("L2").Select()
